$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 12500888
$ws.Range("I9").Value = 50000450
$ws.Range("J9").Value = 1033.3334
$ws.Range("K9").Value = 50000450
$ws.Range("L9").Value = 1033.3334
$ws.Range("M9").Value = -50000281
$ws.Range("N9").Value = -1371.3334
$ws.Range("H62").Value = 3049.3333
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 4111
$ws.Range("K62").Value = 2200
$ws.Range("L62").Value = 4111
$ws.Range("M62").Value = -1576
$ws.Range("N62").Value = -5359
$ws.Range("H65").Value = 3049.3333
$ws.Range("I65").Value = 2200
$ws.Range("J65").Value = 4111
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 20555
$ws.Range("M65").Value = -7880
$ws.Range("N65").Value = -26795
$ws.Range("H99").Value = 1025
$ws.Range("I99").Value = 341.125
$ws.Range("J99").Value = 1708.875
$ws.Range("K99").Value = 1023.375
$ws.Range("L99").Value = 5126.625
$ws.Range("M99").Value = 474.625
$ws.Range("N99").Value = -8122.625
$ws.Range("H100").Value = 4910.0415
$ws.Range("I100").Value = 4564.6
$ws.Range("J100").Value = 5485.778
$ws.Range("K100").Value = 4564.6
$ws.Range("L100").Value = 5485.778
$ws.Range("M100").Value = -4023.6
$ws.Range("N100").Value = -6567.778
$ws.Range("H135").Value = 1029.85
$ws.Range("I135").Value = 949.82355
$ws.Range("K135").Value = 8548.41195
$ws.Range("M135").Value = -6013.41195
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5627.024
$ws.Range("I32").Value = 4222
$ws.Range("J32").Value = 13300.615
$ws.Range("K32").Value = 4222
$ws.Range("L32").Value = 13300.615
$ws.Range("M32").Value = -3935
$ws.Range("N32").Value = -13874.615
$ws.Range("H61").Value = 2789.5217
$ws.Range("I61").Value = 1616.8695
$ws.Range("J61").Value = 3962.1738
$ws.Range("K61").Value = 1616.8695
$ws.Range("L61").Value = 3962.1738
$ws.Range("M61").Value = -1404.8695
$ws.Range("N61").Value = -4386.1738
$ws.Range("H102").Value = 3150.6428
$ws.Range("I102").Value = 2650.75
$ws.Range("J102").Value = 6150
$ws.Range("K102").Value = 2650.75
$ws.Range("L102").Value = 6150
$ws.Range("M102").Value = -1028.75
$ws.Range("N102").Value = -9394
$ws.Range("H122").Value = 2757.6875
$ws.Range("I122").Value = 2078.647
$ws.Range("J122").Value = 3527.2666
$ws.Range("K122").Value = 6235.941
$ws.Range("L122").Value = 10581.7998
$ws.Range("M122").Value = -3785.941
$ws.Range("N122").Value = -15481.7998
$ws.Range("H136").Value = 2789.5217
$ws.Range("I136").Value = 1616.8695
$ws.Range("J136").Value = 3962.1738
$ws.Range("K136").Value = 4850.6085
$ws.Range("L136").Value = 11886.5214
$ws.Range("M136").Value = -2300.6085
$ws.Range("N136").Value = -16986.5214
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 21040.04
$ws.Range("I86").Value = 1075.125
$ws.Range("J86").Value = 52983.9
$ws.Range("K86").Value = 1075.125
$ws.Range("L86").Value = 52983.9
$ws.Range("M86").Value = 47.875
$ws.Range("N86").Value = -55229.9
$ws.Range("H89").Value = 21040.04
$ws.Range("I89").Value = 1075.125
$ws.Range("J89").Value = 52983.9
$ws.Range("K89").Value = 5375.625
$ws.Range("L89").Value = 264919.5
$ws.Range("M89").Value = 240.375
$ws.Range("N89").Value = -276151.5
$ws.Range("H99").Value = 1978.909
$ws.Range("I99").Value = 921
$ws.Range("K99").Value = 921
$ws.Range("M99").Value = 577
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3362.5
$ws.Range("I99").Value = 900
$ws.Range("K99").Value = 900
$ws.Range("M99").Value = 598
$ws.Range("H126").Value = 3362.5
$ws.Range("I126").Value = 900
$ws.Range("K126").Value = 2700
$ws.Range("M126").Value = -230
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 815.4545000000001
$ws.Range("J86").Value = 876
$ws.Range("L86").Value = 2628
$ws.Range("N86").Value = -5000
$ws.Range("H89").Value = 815.4545000000001
$ws.Range("J89").Value = 876
$ws.Range("L89").Value = 7884
$ws.Range("N89").Value = -19740
$ws.Range("H92").Value = 1350
$ws.Range("J92").Value = 1655.5555
$ws.Range("L92").Value = 4966.666499999999
$ws.Range("N92").Value = -7462.666499999999
$ws.Range("H127").Value = 2283.25
$ws.Range("J127").Value = 2283.25
$ws.Range("L127").Value = 6849.75
$ws.Range("N127").Value = -16769.75
$ws.Range("H131").Value = 1222.5555
$ws.Range("I131").Value = 3530.5
$ws.Range("J131").Value = 1037.92
$ws.Range("K131").Value = 10591.5
$ws.Range("L131").Value = 3113.76
$ws.Range("M131").Value = -5551.5
$ws.Range("N131").Value = -13193.76
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4813.3335
$ws.Range("I70").Value = 5052.857
$ws.Range("J70").Value = 3975
$ws.Range("K70").Value = 5052.857
$ws.Range("L70").Value = 3975
$ws.Range("M70").Value = -4782.857
$ws.Range("N70").Value = -4515
$ws.Range("H73").Value = 4813.3335
$ws.Range("I73").Value = 5052.857
$ws.Range("J73").Value = 3975
$ws.Range("K73").Value = 5052.857
$ws.Range("L73").Value = 3975
$ws.Range("M73").Value = -4116.857
$ws.Range("N73").Value = -5847
$ws.Range("H102").Value = 61515.41
$ws.Range("I102").Value = 2157.8462
$ws.Range("J102").Value = 254427.5
$ws.Range("K102").Value = 2157.8462
$ws.Range("L102").Value = 254427.5
$ws.Range("M102").Value = -535.8462
$ws.Range("N102").Value = -257671.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3112
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 3434.4
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 3434.4
$ws.Range("M40").Value = -1364
$ws.Range("N40").Value = -3706.4
$ws.Range("H46").Value = 1287.8572
$ws.Range("J46").Value = 1818.75
$ws.Range("L46").Value = 1818.75
$ws.Range("N46").Value = -2194.75
$ws.Range("H69").Value = 35000
$ws.Range("J69").Value = 35000
$ws.Range("L69").Value = 35000
$ws.Range("N69").Value = -36622
$ws.Range("H72").Value = 35000
$ws.Range("J72").Value = 35000
$ws.Range("L72").Value = 105000
$ws.Range("N72").Value = -113112
$ws.Range("H132").Value = 2726.3
$ws.Range("I132").Value = 1752.25
$ws.Range("J132").Value = 3700.35
$ws.Range("K132").Value = 5256.75
$ws.Range("L132").Value = 11101.05
$ws.Range("M132").Value = -2726.75
$ws.Range("N132").Value = -16161.05
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1944.5
$ws.Range("I96").Value = 1944.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1944.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = ""
$ws.Range("N96").Value = -571.5
$ws.Range("H100").Value = 1005.5
$ws.Range("I100").Value = 877.7143
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 1755.4286
$ws.Range("L100").Value = 3800
$ws.Range("M100").Value = -1214.4286
$ws.Range("N100").Value = -4882
$ws.Range("H122").Value = 716461.2
$ws.Range("I122").Value = 1112901.4
$ws.Range("K122").Value = 3338704.2
$ws.Range("M122").Value = -3336254.2
$ws.Range("H132").Value = 236549.69
$ws.Range("I132").Value = 325859.12
$ws.Range("J132").Value = 38793.07
$ws.Range("K132").Value = 977577.36
$ws.Range("L132").Value = 116379.21
$ws.Range("M132").Value = -975047.36
$ws.Range("N132").Value = -121439.21
